$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column G (pushes old G,H,I -> H,I,J)
$ws.Range("G:G").Insert()

# New header for inserted column
$ws.Range("G2").Value = "dt"

# New dt values for existing rows 3 and 4
$ws.Range("G3").Value = 0.2
$ws.Range("G4").Value = 0.2

# New row 5
$ws.Range("C5").Value = 20
$ws.Range("D5").Value = 10
$ws.Range("E5").Value = 20
$ws.Range("F5").Value = 50
$ws.Range("G5").Value = 0.4
$ws.Range("I5").Value = "too few to join boids stuck"

# New row 6
$ws.Range("C6").Value = 20
$ws.Range("D6").Value = 10
$ws.Range("E6").Value = 20
$ws.Range("F6").Value = 100
$ws.Range("G6").Value = 0.2
$ws.Range("H6").Value = 107
$ws.Range("I6").Value = "enough to run and the boids merge"

# Update selection to match new state
$ws.Range("C7").Select()
